# Actualiza horarios Linea 141 - 1273
# Actualiza encabezados (Ultima actualizacion / Total filas) y las filas de datos
# en las 3 hojas del libro: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

function Set-ScheduleSheet {
    param($ws, $lastUpdate, $totalRows, $rows)

    # Fila 2: Última actualización ; Fila 3: Total filas
    $ws.Cells.Item(2, 1).Value = "Última actualización: " + $lastUpdate
    $ws.Cells.Item(3, 1).Value = "Total filas: " + $totalRows

    # Filas de datos a partir de la fila 6 (fila 5 = encabezados)
    $r = 6
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $r = $r + 1
    }
}

# ---- Hoja LP1912 ----
$wsLP1912 = $wb.Worksheets.Item("LP1912")
$rowsLP1912 = @(
    @("05:47:29","05:47","17_ROMERO",0,"LP1912"),
    @("05:47:29","06:09","10_OLMOS",22,"LP1912"),
    @("06:15:23","06:16","215A_EL PATO",1,"LP1912"),
    @("06:15:23","06:30","23_HERNANDEZ",15,"LP1912"),
    @("06:15:23","06:34","11_ETCHEVERRY",19,"LP1912"),
    @("06:15:23","06:39","17X38_ROMERO",24,"LP1912"),
    @("06:15:23","06:41","16_SANTA ANA",26,"LP1912"),
    @("06:46:40","06:56","215A_EL PATO",10,"LP1912"),
    @("06:15:23","06:57","215A_EL PATO",42,"LP1912"),
    @("06:46:40","06:59","225_GOMEZ",13,"LP1912"),
    @("06:46:40","07:15","215C_EL PATO",29,"LP1912"),
    @("06:15:23","07:16","215C_EL PATO",61,"LP1912"),
    @("06:46:40","07:19","14_ABASTO",33,"LP1912"),
    @("06:46:40","07:20","16_SANTA ANA",34,"LP1912"),
    @("06:15:23","07:21","16_SANTA ANA",66,"LP1912"),
    @("06:46:40","07:21","23_HERNANDEZ",35,"LP1912"),
    @("06:46:40","07:29","17X38_ROMERO",43,"LP1912"),
    @("06:46:40","07:35","10_OLMOS",49,"LP1912"),
    @("06:46:40","07:36","27_EL RETIRO",50,"LP1912"),
    @("06:15:23","07:37","27_EL RETIRO",82,"LP1912"),
    @("06:46:40","07:43","215A_EL PATO",57,"LP1912"),
    @("06:46:40","07:55","14_ABASTO",69,"LP1912"),
    @("06:46:40","08:00","17_ROMERO",74,"LP1912"),
    @("06:46:40","08:01","16_SANTA ANA",75,"LP1912"),
    @("06:46:40","08:06","23_HERNANDEZ",80,"LP1912"),
    @("06:46:40","08:11","10_OLMOS",85,"LP1912"),
    @("06:46:40","08:13","15X38_ABASTO",87,"LP1912"),
    @("06:46:40","08:29","15_ABASTO",103,"LP1912"),
    @("06:46:40","08:29","11_ETCHEVERRY",103,"LP1912"),
    @("06:46:40","08:41","16_P MOR-SANTA ANA",115,"LP1912"),
    @("06:46:40","08:43","215C_EL PATO",117,"LP1912")
)
Set-ScheduleSheet $wsLP1912 "06:46:40" 31 $rowsLP1912

# ---- Hoja LP1912-215 ----
$wsLP1912215 = $wb.Worksheets.Item("LP1912-215")
$rowsLP1912215 = @(
    @("06:15:23","06:16","215A_EL PATO",1,"LP1912"),
    @("06:46:40","06:56","215A_EL PATO",10,"LP1912"),
    @("06:15:23","06:57","215A_EL PATO",42,"LP1912"),
    @("06:46:40","07:15","215C_EL PATO",29,"LP1912"),
    @("06:15:23","07:16","215C_EL PATO",61,"LP1912"),
    @("06:46:40","07:43","215A_EL PATO",57,"LP1912"),
    @("06:46:40","08:43","215C_EL PATO",117,"LP1912")
)
Set-ScheduleSheet $wsLP1912215 "06:46:40" 7 $rowsLP1912215

# ---- Hoja 6203-6173 ----
$ws62036173 = $wb.Worksheets.Item("6203-6173")
$rows62036173 = @(
    @("06:46:40","07:42","215A_LA PLATA",56,"L6173"),
    @("06:15:23","07:43","215A_LA PLATA",88,"L6173"),
    @("06:46:40","08:35","215A_LA PLATA",109,"L6173")
)
Set-ScheduleSheet $ws62036173 "06:46:40" 3 $rows62036173
